$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 10.100659
$ws.Cells.Item(2, 8).Value = 30.301977
$ws.Cells.Item(2, 9).Value = 0.3328245842863797
$ws.Cells.Item(2, 10).Value = 0.3328245842863797
$ws.Cells.Item(2, 13).Value = 16.27546433333333
$ws.Cells.Item(2, 14).Value = 48.826393
$ws.Cells.Item(2, 15).Value = 0.06628560529319844
$ws.Cells.Item(2, 16).Value = 0.06628560529319844
$ws.Cells.Item(2, 17).Value = 164.3929152976623
$ws.Cells.Item(2, 18).Value = 1479.536237678961
$ws.Cells.Item(2, 19).Value = 0.02206147902587982
$ws.Cells.Item(2, 20).Value = 0.02206147902587982
# Row 3
$ws.Cells.Item(3, 7).Value = 10.100659
$ws.Cells.Item(3, 8).Value = 30.301977
$ws.Cells.Item(3, 9).Value = 0.3328245842863797
$ws.Cells.Item(3, 10).Value = 0.3328245842863797
$ws.Cells.Item(3, 15).Value = 0.3480686258826592
$ws.Cells.Item(3, 16).Value = 0.3480686258826592
$ws.Cells.Item(3, 17).Value = 863.2344214011885
$ws.Cells.Item(3, 18).Value = 7769.109792610696
$ws.Cells.Item(3, 19).Value = 0.1158457957125275
$ws.Cells.Item(3, 20).Value = 0.1158457957125275
# Row 4
$ws.Cells.Item(4, 7).Value = 10.100659
$ws.Cells.Item(4, 8).Value = 30.301977
$ws.Cells.Item(4, 9).Value = 0.3328245842863797
$ws.Cells.Item(4, 10).Value = 0.3328245842863797
$ws.Cells.Item(4, 13).Value = 42.61351133333333
$ws.Cells.Item(4, 14).Value = 127.840534
$ws.Cells.Item(4, 15).Value = 0.17355341356458
$ws.Cells.Item(4, 16).Value = 0.17355341356458
$ws.Cells.Item(4, 17).Value = 430.4245467706353
$ws.Cells.Item(4, 18).Value = 3873.820920935718
$ws.Cells.Item(4, 19).Value = 0.05776284272111346
$ws.Cells.Item(4, 20).Value = 0.05776284272111347
# Row 5
$ws.Cells.Item(5, 7).Value = 10.100659
$ws.Cells.Item(5, 8).Value = 30.301977
$ws.Cells.Item(5, 9).Value = 0.3328245842863797
$ws.Cells.Item(5, 10).Value = 0.3328245842863797
$ws.Cells.Item(5, 13).Value = 101.183272
$ws.Cells.Item(5, 14).Value = 303.549816
$ws.Cells.Item(5, 15).Value = 0.4120923552595624
$ws.Cells.Item(5, 16).Value = 0.4120923552595624
$ws.Cells.Item(5, 17).Value = 1022.017726976248
$ws.Cells.Item(5, 18).Value = 9198.15954278623
$ws.Cells.Item(5, 19).Value = 0.1371544668268589
$ws.Cells.Item(5, 20).Value = 0.1371544668268589
# Row 6
$ws.Cells.Item(6, 9).Value = 0.4180918757349671
$ws.Cells.Item(6, 10).Value = 0.4180918757349671
$ws.Cells.Item(6, 13).Value = 16.27546433333333
$ws.Cells.Item(6, 14).Value = 48.826393
$ws.Cells.Item(6, 15).Value = 0.06628560529319844
$ws.Cells.Item(6, 16).Value = 0.06628560529319844
$ws.Cells.Item(6, 17).Value = 206.5092110359226
$ws.Cells.Item(6, 18).Value = 1858.582899323304
$ws.Cells.Item(6, 19).Value = 0.027713473051261
$ws.Cells.Item(6, 20).Value = 0.027713473051261
# Row 7
$ws.Cells.Item(7, 9).Value = 0.4180918757349671
$ws.Cells.Item(7, 10).Value = 0.4180918757349671
$ws.Cells.Item(7, 15).Value = 0.3480686258826592
$ws.Cells.Item(7, 16).Value = 0.3480686258826592
$ws.Cells.Item(7, 19).Value = 0.1455246646797735
$ws.Cells.Item(7, 20).Value = 0.1455246646797735
# Row 8
$ws.Cells.Item(8, 9).Value = 0.4180918757349671
$ws.Cells.Item(8, 10).Value = 0.4180918757349671
$ws.Cells.Item(8, 13).Value = 42.61351133333333
$ws.Cells.Item(8, 14).Value = 127.840534
$ws.Cells.Item(8, 15).Value = 0.17355341356458
$ws.Cells.Item(8, 16).Value = 0.17355341356458
$ws.Cells.Item(8, 17).Value = 540.6962544775946
$ws.Cells.Item(8, 18).Value = 4866.266290298352
$ws.Cells.Item(8, 19).Value = 0.07256127221742173
$ws.Cells.Item(8, 20).Value = 0.07256127221742174
# Row 9
$ws.Cells.Item(9, 9).Value = 0.4180918757349671
$ws.Cells.Item(9, 10).Value = 0.4180918757349671
$ws.Cells.Item(9, 13).Value = 101.183272
$ws.Cells.Item(9, 14).Value = 303.549816
$ws.Cells.Item(9, 15).Value = 0.4120923552595624
$ws.Cells.Item(9, 16).Value = 0.4120923552595624
$ws.Cells.Item(9, 17).Value = 1283.851400046272
$ws.Cells.Item(9, 18).Value = 11554.66260041645
$ws.Cells.Item(9, 19).Value = 0.1722924657865109
$ws.Cells.Item(9, 20).Value = 0.1722924657865109
# Row 10
$ws.Cells.Item(10, 7).Value = 4.721016333333334
$ws.Cells.Item(10, 8).Value = 14.163049
$ws.Cells.Item(10, 9).Value = 0.1555611667071302
$ws.Cells.Item(10, 10).Value = 0.1555611667071302
$ws.Cells.Item(10, 13).Value = 16.27546433333333
$ws.Cells.Item(10, 14).Value = 48.826393
$ws.Cells.Item(10, 15).Value = 0.06628560529319844
$ws.Cells.Item(10, 16).Value = 0.06628560529319844
$ws.Cells.Item(10, 17).Value = 76.83673295025078
$ws.Cells.Item(10, 18).Value = 691.530596552257
$ws.Cells.Item(10, 19).Value = 0.01031146609529828
$ws.Cells.Item(10, 20).Value = 0.01031146609529827
# Row 11
$ws.Cells.Item(11, 7).Value = 4.721016333333334
$ws.Cells.Item(11, 8).Value = 14.163049
$ws.Cells.Item(11, 9).Value = 0.1555611667071302
$ws.Cells.Item(11, 10).Value = 0.1555611667071302
$ws.Cells.Item(11, 15).Value = 0.3480686258826592
$ws.Cells.Item(11, 16).Value = 0.3480686258826592
$ws.Cells.Item(11, 17).Value = 403.4730608102462
$ws.Cells.Item(11, 18).Value = 3631.257547292215
$ws.Cells.Item(11, 19).Value = 0.0541459615364541
$ws.Cells.Item(11, 20).Value = 0.05414596153645408
# Row 12
$ws.Cells.Item(12, 7).Value = 4.721016333333334
$ws.Cells.Item(12, 8).Value = 14.163049
$ws.Cells.Item(12, 9).Value = 0.1555611667071302
$ws.Cells.Item(12, 10).Value = 0.1555611667071302
$ws.Cells.Item(12, 13).Value = 42.61351133333333
$ws.Cells.Item(12, 14).Value = 127.840534
$ws.Cells.Item(12, 15).Value = 0.17355341356458
$ws.Cells.Item(12, 16).Value = 0.17355341356458
$ws.Cells.Item(12, 17).Value = 201.1790830253518
$ws.Cells.Item(12, 18).Value = 1810.611747228166
$ws.Cells.Item(12, 19).Value = 0.02699817150011114
$ws.Cells.Item(12, 20).Value = 0.02699817150011114
# Row 13
$ws.Cells.Item(13, 7).Value = 4.721016333333334
$ws.Cells.Item(13, 8).Value = 14.163049
$ws.Cells.Item(13, 9).Value = 0.1555611667071302
$ws.Cells.Item(13, 10).Value = 0.1555611667071302
$ws.Cells.Item(13, 13).Value = 101.183272
$ws.Cells.Item(13, 14).Value = 303.549816
$ws.Cells.Item(13, 15).Value = 0.4120923552595624
$ws.Cells.Item(13, 16).Value = 0.4120923552595624
$ws.Cells.Item(13, 17).Value = 477.6878797721093
$ws.Cells.Item(13, 18).Value = 4299.190917948984
$ws.Cells.Item(13, 19).Value = 0.06410556757526673
$ws.Cells.Item(13, 20).Value = 0.06410556757526671
# Row 14
$ws.Cells.Item(14, 7).Value = 2.838244666666667
$ws.Cells.Item(14, 8).Value = 8.514734000000001
$ws.Cells.Item(14, 9).Value = 0.09352237327152295
$ws.Cells.Item(14, 10).Value = 0.09352237327152294
$ws.Cells.Item(14, 13).Value = 16.27546433333333
$ws.Cells.Item(14, 14).Value = 48.826393
$ws.Cells.Item(14, 15).Value = 0.06628560529319844
$ws.Cells.Item(14, 16).Value = 0.06628560529319844
$ws.Cells.Item(14, 17).Value = 46.19374984160689
$ws.Cells.Item(14, 18).Value = 415.743748574462
$ws.Cells.Item(14, 19).Value = 0.006199187120759342
$ws.Cells.Item(14, 20).Value = 0.00619918712075934
# Row 15
$ws.Cells.Item(15, 7).Value = 2.838244666666667
$ws.Cells.Item(15, 8).Value = 8.514734000000001
$ws.Cells.Item(15, 9).Value = 0.09352237327152295
$ws.Cells.Item(15, 10).Value = 0.09352237327152294
$ws.Cells.Item(15, 15).Value = 0.3480686258826592
$ws.Cells.Item(15, 16).Value = 0.3480686258826592
$ws.Cells.Item(15, 17).Value = 242.5654101009656
$ws.Cells.Item(15, 18).Value = 2183.08869090869
$ws.Cells.Item(15, 19).Value = 0.03255220395390412
$ws.Cells.Item(15, 20).Value = 0.03255220395390412
# Row 16
$ws.Cells.Item(16, 7).Value = 2.838244666666667
$ws.Cells.Item(16, 8).Value = 8.514734000000001
$ws.Cells.Item(16, 9).Value = 0.09352237327152295
$ws.Cells.Item(16, 10).Value = 0.09352237327152294
$ws.Cells.Item(16, 13).Value = 42.61351133333333
$ws.Cells.Item(16, 14).Value = 127.840534
$ws.Cells.Item(16, 15).Value = 0.17355341356458
$ws.Cells.Item(16, 16).Value = 0.17355341356458
$ws.Cells.Item(16, 17).Value = 120.9475712697729
$ws.Cells.Item(16, 18).Value = 1088.528141427956
$ws.Cells.Item(16, 19).Value = 0.01623112712593364
$ws.Cells.Item(16, 20).Value = 0.01623112712593364
# Row 17
$ws.Cells.Item(17, 7).Value = 2.838244666666667
$ws.Cells.Item(17, 8).Value = 8.514734000000001
$ws.Cells.Item(17, 9).Value = 0.09352237327152295
$ws.Cells.Item(17, 10).Value = 0.09352237327152294
$ws.Cells.Item(17, 13).Value = 101.183272
$ws.Cells.Item(17, 14).Value = 303.549816
$ws.Cells.Item(17, 15).Value = 0.4120923552595624
$ws.Cells.Item(17, 16).Value = 0.4120923552595624
$ws.Cells.Item(17, 17).Value = 287.1828821098827
$ws.Cells.Item(17, 18).Value = 2584.645938988944
$ws.Cells.Item(17, 19).Value = 0.03853985507092582
$ws.Cells.Item(17, 20).Value = 0.03853985507092583
